$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add the new BOM row (item + URL) at the end of the table
$newRow = $lo.ListRows.Add()
$newRow.Range.Cells.Item(1,1).Value = "Rocker switch"
$newRow.Range.Cells.Item(1,2).Value = "https://smile.amazon.com/gp/product/B07D1J246N/ref=oh_aui_search_detailpage?ie=UTF8&psc=1"

# Re-sort the table by Item (column A), ascending, keeping the header row in place
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A21"))
$sortObj.SetRange($ws.Range("A1:B21"))
$sortObj.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$sortObj.Apply()

$ws.Range("B17").Select() | Out-Null
